# feat: Add delivery package
#
# - Rename the existing product "Chapelet" to
#   "Attâche câbles à crochet et à boucle" (applies to every existing order
#   row that references it, rows 2-7, and fills in the previously-empty
#   E7 product cell).
# - Append three new orders (rows 8-10) for HARRY POTTER, SATORU GOJO and
#   JOCODE DEV.
# - Widen column E to fit the new product name.
# - Select the full data range, mirroring the author's final selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newProduct = "Attâche câbles à crochet et à boucle"
$phone = "2250788466748"
$address = "Koumassi Quartier Divo"
$addressLong = "Koumassi Quartier Divo carrefour canniveau - premier carrefour à gauche"

# --- Rename the "Chapelet" product everywhere it is used (rows 2-7) ----
# All of E2:E6 already hold the old product name, and E7 was previously
# blank; setting the whole block at once keeps every row pointing at the
# same (renamed) shared string instead of creating a duplicate entry.
$ws.Range("E2:E7").Value = $newProduct

# --- Row 8: HARRY POTTER -------------------------------------------------
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "HARRY POTTER"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = $phone
$ws.Range("D8").Value = $address
$ws.Range("E8").Value = $newProduct
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 5000

# --- Row 9: SATORU GOJO ---------------------------------------------------
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "SATORU GOJO"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = $phone
$ws.Range("D9").Value = $address
$ws.Range("E9").Value = $newProduct
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5000

# --- Row 10: JOCODE DEV ---------------------------------------------------
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "JOCODE DEV"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = $phone
$ws.Range("D10").Value = $address
$ws.Range("E10").Value = $newProduct
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5000

# Row 9's address gets refined afterwards to the longer, more precise
# description (done last so the new shared strings land in the same
# order the original authoring session produced them in).
$ws.Range("D9").Value = $addressLong

# --- Cosmetics -------------------------------------------------------------
# Fit column E to the new product name's width.
$ws.Columns.Item(5).ColumnWidth = 30.14

# Select the whole updated table, as the author's session left it.
$ws.Range("A1:G10").Select() | Out-Null
